$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lookup table data (Table2: G1:H4) -----------------------------------
# Header row (G1/H1) already contains "floor type " / "cost/sqft".
$ws.Range("G2").Value = "hardwood"
$ws.Range("H2").Value = 1.39
$ws.Range("G3").Value = "carpet"
$ws.Range("H3").Value = 3.99
$ws.Range("G4").Value = "tile "
$ws.Range("H4").Value = 4.99

# --- Test cases (rows 3-7) ------------------------------------------------
# Row 3: hardwood, 1 x 1
$ws.Range("A3").Formula = "=H3"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1

# Row 4: tile, 5 x 5
$ws.Range("A4").Formula = "=H4"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 5

# Row 5: tile, 6 x 8
$ws.Range("A5").Formula = "=H4"
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 8

# Row 6: hardwood, 4 x 10
$ws.Range("A6").Formula = "=H3"
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 10

# Row 7: carpet, 2 x 4
$ws.Range("A7").Formula = "=H2"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 4

# D4:D7 share one formula group (room cost = price * length * width)
$ws.Range("D4:D7").Formula = "=A4*B4*C4"
# D3 has its own (non-shared) instance of the same formula
$ws.Range("D3").Formula = "=A3*B3*C3"

# Row 8: final overall cost
$ws.Range("D8").Formula = "=D3+D4+D5+D6+D7"

# --- Misc view state -------------------------------------------------------
$ws.Range("C8").Select()
